{"js": "// Update the Notetaker role paragraph: make the \"choose a side\" rule\n// conditional on team size (3+ members vs. exactly 2 members).\n\nconst body = context.document.body;\n\n// The run we need to replace contains this exact (unique) sentence.\nconst oldText =\n  \"The Notetaker must also choose either the Employee or Employer side. \" +\n  \"Write the Notetaker\\u2019s name in both the Notetaker row AND their chosen advocate row.\";\n\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the Notetaker paragraph text to replace.\");\n}\n\nconst target = results.items[0];\n\n// The new text is split across several runs so that the \"3 or more members\"\n// and \"2 members\" phrases can be bolded, matching the rest of the document's\n// existing bold-callout style (e.g. \"Tip:\", \"Note:\").\nconst segments = [\n  { text: \"If your team has\", bold: false },\n  { text: \" \", bold: false },\n  { text: \"3 or more members\", bold: true },\n  { text: \" \", bold: false },\n  {\n    text:\n      \"present, the Notetaker focuses only on facilitating and recording the debate. If only\",\n    bold: false,\n  },\n  { text: \" \", bold: false },\n  { text: \"2 members\", bold: true },\n  { text: \" \", bold: false },\n  {\n    text:\n      \"are present, the Notetaker must also choose either the Employee or Employer side. \" +\n      \"In that case, write the Notetaker\\u2019s name in both the Notetaker row AND their chosen advocate row.\",\n    bold: false,\n  },\n];\n\n// Replace the whole old run with the first segment, then insert each\n// subsequent segment right after the previous one, setting bold explicitly\n// on every piece so formatting never leaks between segments.\nlet current = target.insertText(segments[0].text, Word.InsertLocation.replace);\ncurrent.font.bold = segments[0].bold;\nawait context.sync();\n\nfor (let i = 1; i < segments.length; i++) {\n  current = current.insertText(segments[i].text, Word.InsertLocation.after);\n  current.font.bold = segments[i].bold;\n  await context.sync();\n}\n", "ps1": "# Update the Notetaker role paragraph: make the \"choose a side\" rule\n# conditional on team size (3+ members vs. exactly 2 members).\n\n$d = $word.ActiveDocument\n\n$rightQuote = [char]0x2019\n\n$oldText = \"The Notetaker must also choose either the Employee or Employer side. \" `\n    + \"Write the Notetaker\" + $rightQuote + \"s name in both the Notetaker row AND their chosen advocate row.\"\n\n$rng = $d.Content\n$found = $rng.Find.Execute($oldText)\nif (-not $found) {\n    throw \"Could not find the Notetaker paragraph text to replace.\"\n}\n\n# The new text is split across several runs so that the \"3 or more members\"\n# and \"2 members\" phrases can be bolded, matching the rest of the document's\n# existing bold-callout style (e.g. \"Tip:\", \"Note:\").\n$rng.Text = \"If your team has\"\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\" \")\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\"3 or more members\")\n$rng.Bold = 1\n$rng.Collapse(0)\n\n$rng.InsertAfter(\" \")\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\"present, the Notetaker focuses only on facilitating and recording the debate. If only\")\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\" \")\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\"2 members\")\n$rng.Bold = 1\n$rng.Collapse(0)\n\n$rng.InsertAfter(\" \")\n$rng.Bold = 0\n$rng.Collapse(0)\n\n$rng.InsertAfter(\"are present, the Notetaker must also choose either the Employee or Employer side. \" `\n    + \"In that case, write the Notetaker\" + $rightQuote + \"s name in both the Notetaker row AND their chosen advocate row.\")\n$rng.Bold = 0\n"}
